$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.429.66"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "2.253.56"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'308.00"
$ws.Range("E5").Value = "  -1.93%  "
$ws.Range("D6").Value = "97.31"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D7").Value = "0.575"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.530"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").Value = "35.31"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").Value = "0.0817"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "7.32"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").Value = "2.596.62"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "2.252.83"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").Value = "0.839"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "13.69"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").Value = "44.222.88"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").Value = "0.0₃0973"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "6.42"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "12.23"
$ws.Range("E21").Value = "  -6.34%  "
$ws.Range("D22").Value = "65.84"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "238.52"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "38.99"
$ws.Range("E27").Value = "  +6.73%  "
$ws.Range("D28").Value = "10.01"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  +3.64%  "
$ws.Range("D30").Value = "5.97"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "20.18"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("D32").Value = "152.95"
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("D33").Value = "0.0802"
$ws.Range("E33").Value = "  -3.55%  "
$ws.Range("D34").Value = "3.25"
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D35").Value = "2.62"
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("E37").Value = "  -1.69%  "
$ws.Range("D38").Value = "1.79"
$ws.Range("E38").Value = "  -5.80%  "
$ws.Range("D39").Value = "3.63"
$ws.Range("E39").Value = "  +2.69%  "
$ws.Range("D40").Value = "14.75"
$ws.Range("E40").Value = "  -6.31%  "
$ws.Range("D41").Value = "3.89"
$ws.Range("E41").Value = "  -3.33%  "
$ws.Range("D42").Value = "0.0302"
$ws.Range("E42").Value = "  -1.46%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "1.756.36"
$ws.Range("E44").Value = "  +3.18%  "
$ws.Range("D45").Value = "83.62"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "100.93"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").Value = "'5.00"
$ws.Range("E48").Value = "  -2.74%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "8.19"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.60"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "55.21"
$ws.Range("E51").Value = "  -1.81%  "
